# "Generate Report for Archive"
# Swaps the report rows for the two localization entries
#   6e6f5a0a-f847-4138-9cfe-7cdb61058920  (row 9 after the edit)
#   e7c36651-26d0-49ee-b1cf-7f40f5dd3b38  (row 8 after the edit)
# across the Overview / zh-cn / de-de sheets, and refreshes the status
# of e7c36651 from "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview": columns A File Name, B Path And Name (hyperlink),
# C Extension, D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date
# Only A, B, E, F, G actually change between the two rows - C/D are
# identical for both entries so they are left untouched.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A8").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
$wsOverview.Range("B8").Value = "e2e\e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
$wsOverview.Range("E8").Value = "In Translation"
$wsOverview.Range("F8").Value = "In Translation"
$wsOverview.Range("G8").Value = "2016-12-16 08:22:28"

$wsOverview.Range("A9").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
$wsOverview.Range("B9").Value = "e2e\6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-12-16 08:12:14"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Row -eq 8) {
        $hl.TextToDisplay = "e2e\e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
    } elseif ($hl.Range.Row -eq 9) {
        $hl.TextToDisplay = "e2e\6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
    }
}

# ---------------------------------------------------------------
# Sheet "zh-cn": columns A Source File Name, C Status,
# G Latest Handoff File, H Latest Handoff Datetime change;
# the rest of the row is identical between the two entries.
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A8").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
$wsZhCn.Range("C8").Value = "In Translation"
$wsZhCn.Range("G8").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.ebe39ec3a11a72ce0470bafa1cc822a30a67b978.zh-cn.xlf"
$wsZhCn.Range("H8").Value = "2016-12-16 08:22:14"

$wsZhCn.Range("A9").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("G9").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.fa7b938b8ca0282e071b9dfae621037cafe4c44e.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-12-16 08:12:01"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Row -eq 8) {
        $hl.TextToDisplay = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
    } elseif ($hl.Range.Row -eq 9) {
        $hl.TextToDisplay = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
    }
}

# ---------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn.
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A8").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
$wsDeDe.Range("C8").Value = "In Translation"
$wsDeDe.Range("G8").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.ebe39ec3a11a72ce0470bafa1cc822a30a67b978.de-de.xlf"
$wsDeDe.Range("H8").Value = "2016-12-16 08:22:28"

$wsDeDe.Range("A9").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("G9").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.fa7b938b8ca0282e071b9dfae621037cafe4c44e.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-12-16 08:12:14"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Row -eq 8) {
        $hl.TextToDisplay = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
    } elseif ($hl.Range.Row -eq 9) {
        $hl.TextToDisplay = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
    }
}
